# "arrumando os champs certo"
#
# 1) Column F (ASSISTS) on rows 2-41 was stored as text ("0", "1", "2", ...)
#    instead of a real number. Re-write each cell with its numeric value so
#    Excel stores it as a number (t="n") rather than an inline string.
# 2) Several rows in column H (CHAMPION) were tagged with the wrong
#    champion name; they should all read "Vi".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> correct numeric ASSISTS value
$assists = @{
    2 = 0;  3 = 0;  4 = 0;  5 = 0;  6 = 0;  7 = 0;  8 = 0;  9 = 0;  10 = 0
    11 = 1; 12 = 1
    13 = 2; 14 = 2; 15 = 2; 16 = 2; 17 = 2; 18 = 2; 19 = 2; 20 = 2
    21 = 4; 22 = 4; 23 = 4; 24 = 4
    25 = 5; 26 = 5; 27 = 5; 28 = 5; 29 = 5; 30 = 5; 31 = 5; 32 = 5; 33 = 5; 34 = 5; 35 = 5
    36 = 6; 37 = 6; 38 = 6; 39 = 6; 40 = 6; 41 = 6
}

foreach ($row in $assists.Keys) {
    $ws.Cells.Item($row, 6).Value = $assists[$row]
}

# rows whose CHAMPION (column H) should be corrected to "Vi"
$championFixRows = @(5, 11, 17, 23, 28, 29, 30, 31, 35, 41)

foreach ($row in $championFixRows) {
    $ws.Cells.Item($row, 8).Value = "Vi"
}
